$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.951.50'
$ws.Range('E2').Value = '  -0.86%  '
$ws.Range('D3').Value = '1.631.50'
$ws.Range('E3').Value = '  -2.57%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '209.69'
$ws.Range('E5').Value = '  -0.90%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5210'
$ws.Range('E6').Value = '  -0.63%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.002'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2568'
$ws.Range('E8').Value = '  -3.14%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06225'
$ws.Range('E9').Value = '  -0.97%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.23'
$ws.Range('E10').Value = '  -5.03%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07559'
$ws.Range('E11').Value = '  -0.03%  '
$ws.Range('D12').Value = '1.630.38'
$ws.Range('E12').Value = '  -2.46%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.352'
$ws.Range('E13').Value = '  -2.18%  '
$ws.Range('D14').Value = '1.857.59'
$ws.Range('E14').Value = '  -2.41%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5420'
$ws.Range('E15').Value = '  -3.25%  '
$ws.Range('D16').Value = '0.0₅7933'
$ws.Range('E16').Value = '  -0.82%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.65'
$ws.Range('E17').Value = '  -3.21%  '
$ws.Range('D18').Value = '25.937.07'
$ws.Range('E18').Value = '  -1.13%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.003'
$ws.Range('E19').Value = '  -0.05%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.613'
$ws.Range('E20').Value = '  -4.11%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '184.69'
$ws.Range('E21').Value = '  -1.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.02'
$ws.Range('E22').Value = '  -3.81%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.065'
$ws.Range('E23').Value = '  -1.79%  '
$ws.Range('E24').Value = '  -0.19%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.56'
$ws.Range('E25').Value = '  -2.74%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1211'
$ws.Range('E26').Value = '  -2.72%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.342'
$ws.Range('E27').Value = '  -2.85%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.51'
$ws.Range('E28').Value = '  -3.16%  '
$ws.Range('E29').Value = '  +0.58%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05925'
$ws.Range('E30').Value = '  -3.72%  '
$ws.Range('E31').Value = '  -3.36%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.348'
$ws.Range('E32').Value = '  -2.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.352'
$ws.Range('E33').Value = '  -3.93%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.602'
$ws.Range('E34').Value = '  -1.59%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9696'
$ws.Range('E35').Value = '  -2.92%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.719'
$ws.Range('E37').Value = '  -0.63%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5765'
$ws.Range('E38').Value = '  -4.82%  '
$ws.Range('E39').Value = '  -1.18%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.002'
$ws.Range('E40').Value = '  -0.44%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8340'
$ws.Range('E41').Value = '  -4.39%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.598'
$ws.Range('E42').Value = '  -7.84%  '
$ws.Range('D43').Value = '1.008.72'
$ws.Range('E43').Value = '  -6.66%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.56'
$ws.Range('E44').Value = '  -0.37%  '
$ws.Range('D45').Value = '1.782.89'
$ws.Range('E45').Value = '  -2.30%  '
$ws.Range('D46').Value = '0.0₈107'
$ws.Range('E46').Value = '  -1.83%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.9970'
$ws.Range('E47').Value = '  -0.51%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '54.23'
$ws.Range('E48').Value = '  -3.18%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.942'
$ws.Range('E49').Value = '  -0.61%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05170'
$ws.Range('E50').Value = '  -1.23%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4222'
$ws.Range('E51').Value = '  -0.86%  '
